# Auto-generated edit script: updates crypto price/volume columns
# per the commit diff (Sat Mar 30 21:18:13 UTC 2024 GitHub Actions run).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force "number-looking" price strings to stay literal text (preserve
# trailing zeros / exact precision instead of Excel coercing to a float).
$textCells = @("D4","D5","D6","D7","D9","D10","D11","D12","D13","D15","D17","D18","D21","D22","D23","D24","D25","D26","D27","D28","D29","D30","D31","D32","D33","D34","D36","D39","D40","D41","D42","D44","D45","D46","D48","D49","D50")
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range("D2").Value = "69.553.18"
$ws.Range("E2").Value = "  +0.09%  "

$ws.Range("D3").Value = "3.496.12"
$ws.Range("E3").Value = "  +0.06%  "

$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").Value = "602.22"
$ws.Range("E5").Value = "  -1.71%  "

$ws.Range("D6").Value = "194.05"
$ws.Range("E6").Value = "  +2.54%  "

$ws.Range("D7").Value = "0.622"
$ws.Range("E7").Value = "  -0.59%  "

$ws.Range("E8").Value = "  -0.06%  "

$ws.Range("D9").Value = "0.200"
$ws.Range("E9").Value = "  -5.37%  "

$ws.Range("D10").Value = "0.646"
$ws.Range("E10").Value = "  +0.09%  "

$ws.Range("D11").Value = "53.07"
$ws.Range("E11").Value = "  +0.44%  "

$ws.Range("D12").Value = "0.0000299"
$ws.Range("E12").Value = "  -2.39%  "

$ws.Range("D13").Value = "9.46"
$ws.Range("E13").Value = "  -0.09%  "

$ws.Range("D14").Value = "4.063.28"
$ws.Range("E14").Value = "  +0.24%  "

$ws.Range("D15").Value = "591.52"
$ws.Range("E15").Value = "  -1.77%  "

$ws.Range("D16").Value = "69.754.85"
$ws.Range("E16").Value = "  +0.33%  "

$ws.Range("D17").Value = "18.97"
$ws.Range("E17").Value = "  +0.21%  "

$ws.Range("D18").Value = "12.66"
$ws.Range("E18").Value = "  +0.81%  "

$ws.Range("E19").Value = "  +2.31%  "

$ws.Range("D20").Value = "3.490.79"
$ws.Range("E20").Value = "  -0.62%  "

$ws.Range("D21").Value = "0.982"
$ws.Range("E21").Value = "  -0.12%  "

$ws.Range("D22").Value = "18.21"
$ws.Range("E22").Value = "  +6.87%  "

$ws.Range("D23").Value = "5.31"
$ws.Range("E23").Value = "  +3.76%  "

$ws.Range("D24").Value = "4.62"
$ws.Range("E24").Value = "  -1.75%  "

$ws.Range("D25").Value = "101.22"
$ws.Range("E25").Value = "  -4.27%  "

$ws.Range("D26").Value = "3.14"
$ws.Range("E26").Value = "  +3.99%  "

$ws.Range("D27").Value = "10.79"
$ws.Range("E27").Value = "  -1.33%  "

$ws.Range("D28").Value = "9.47"
$ws.Range("E28").Value = "  -1.96%  "

$ws.Range("D29").Value = "33.04"
$ws.Range("E29").Value = "  -0.89%  "

$ws.Range("D30").Value = "4.28"
$ws.Range("E30").Value = "  +5.19%  "

$ws.Range("D31").Value = "7.00"
$ws.Range("E31").Value = "  +1.39%  "

$ws.Range("D32").Value = "12.32"
$ws.Range("E32").Value = "  -1.42%  "

$ws.Range("D33").Value = "0.114"
$ws.Range("E33").Value = "  -0.09%  "

$ws.Range("D34").Value = "63.10"
$ws.Range("E34").Value = "  -0.12%  "

$ws.Range("E38").Value = "  +0.01%  "

$ws.Range("D39").Value = "3.65"
$ws.Range("E39").Value = "  +0.82%  "

$ws.Range("D40").Value = "0.389"
$ws.Range("E40").Value = "  -0.92%  "

$ws.Range("D41").Value = "36.13"
$ws.Range("E41").Value = "  -1.21%  "

$ws.Range("D42").Value = "486.49"
$ws.Range("E42").Value = "  -2.59%  "

$ws.Range("E43").Value = "  -1.80%  "

$ws.Range("D44").Value = "0.0451"
$ws.Range("E44").Value = "  -2.40%  "

$ws.Range("D45").Value = "0.139"
$ws.Range("E45").Value = "  -1.44%  "

$ws.Range("D46").Value = "2.80"
$ws.Range("E46").Value = "  -3.72%  "

$ws.Range("E47").Value = "  -1.22%  "

$ws.Range("D48").Value = "1.01"
$ws.Range("E48").Value = "  +0.35%  "

$ws.Range("D49").Value = "8.38"
$ws.Range("E49").Value = "  -4.17%  "

$ws.Range("D50").Value = "0.000244"
$ws.Range("E50").Value = "  +2.40%  "

$ws.Range("E51").Value = "  +10.14%  "

# Rows 35-37: coin ranking reshuffled this run (Maker now 35th, Fetch.AI 36th,
# PEPE slipped to 37th). Ranks in column A stay put; refresh name/link/price/vol.
$ws.Range("B35").Value = "Maker"
$ws.Range("C35").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D35").Value = "3.721.42"
$ws.Range("E35").Value = "  +3.02%  "

$ws.Range("B36").Value = "Fetch.AI"
$ws.Range("C36").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D36").Value = "3.10"
$ws.Range("E36").Value = "  -1.22%  "

$ws.Range("B37").Value = "PEPE"
$ws.Range("C37").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D37").Value = "0.0₃0813"
$ws.Range("E37").Value = "  +5.42%  "
